# Remove the "Cohort" column from the CasesTab query text in cell B2 of the
# "startup" sheet (the query no longer joins/returns the cohort description).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Irish Terrier','Miniature Pinscher','Poodle']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in ['Bladder, Urethra']`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"
